$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 199.66667
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H4").Value = 1014.4
$ws.Range("I4").Value = 321
$ws.Range("K4").Value = 321
$ws.Range("M4").Value = -207
$ws.Range("H12").Value = 872.375
$ws.Range("I12").Value = 1296.8
$ws.Range("K12").Value = 1296.8
$ws.Range("M12").Value = -1126.8
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H41").Value = 68202.734
$ws.Range("I41").Value = 1373.3334
$ws.Range("J41").Value = 112755.664
$ws.Range("K41").Value = 1373.3334
$ws.Range("L41").Value = 112755.664
$ws.Range("M41").Value = -933.3334
$ws.Range("N41").Value = -113635.664
$ws.Range("H51").Value = 3059.8
$ws.Range("I51").Value = 3149.5
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 3149.5
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -2665.5
$ws.Range("N51").Value = -3968
$ws.Range("H53").Value = 540.9474
$ws.Range("J53").Value = 674.2
$ws.Range("L53").Value = 674.2
$ws.Range("N53").Value = -1948.2
$ws.Range("H70").Value = 9055.299999999999
$ws.Range("I70").Value = 1138.25
$ws.Range("K70").Value = 3414.75
$ws.Range("M70").Value = -3144.75
$ws.Range("H73").Value = 9055.299999999999
$ws.Range("I73").Value = 1138.25
$ws.Range("K73").Value = 3414.75
$ws.Range("M73").Value = -2478.75
$ws.Range("H80").Value = 595.2759
$ws.Range("I80").Value = 463.3846
$ws.Range("K80").Value = 1390.1538
$ws.Range("M80").Value = -392.1538
$ws.Range("H83").Value = 595.2759
$ws.Range("I83").Value = 463.3846
$ws.Range("K83").Value = 4170.4614
$ws.Range("M83").Value = 821.5385999999999
$ws.Range("H97").Value = 499.5
$ws.Range("J97").Value = 599
$ws.Range("L97").Value = 1797
$ws.Range("N97").Value = -2789
$ws.Range("I107").Value = 540.5909
$ws.Range("J107").Value = 288.75
$ws.Range("K107").Value = 540.5909
$ws.Range("L107").Value = 288.75
$ws.Range("M107").Value = 1379.4091
$ws.Range("N107").Value = -4128.75
$ws.Range("H111").Value = 2771.5557
$ws.Range("I111").Value = 358.85715
$ws.Range("J111").Value = 11216
$ws.Range("K111").Value = 1076.57145
$ws.Range("L111").Value = 33648
$ws.Range("M111").Value = 1990.42855
$ws.Range("N111").Value = -39782
$ws.Range("H112").Value = 1331.68
$ws.Range("J112").Value = 1237.7142
$ws.Range("L112").Value = 3713.1426
$ws.Range("N112").Value = -5929.142599999999
$ws.Range("H116").Value = 347761.75
$ws.Range("I116").Value = 94281.336
$ws.Range("J116").Value = 601242.2
$ws.Range("K116").Value = 94281.336
$ws.Range("L116").Value = 601242.2
$ws.Range("M116").Value = -90839.336
$ws.Range("N116").Value = -608126.2
$ws.Range("H125").Value = 580.41174
$ws.Range("I125").Value = 392.3
$ws.Range("K125").Value = 3530.7
$ws.Range("M125").Value = -1070.7
$ws.Range("H132").Value = 34628.855
$ws.Range("I132").Value = 38684.402
$ws.Range("K132").Value = 116053.206
$ws.Range("M132").Value = -113523.206
$ws.Range("H135").Value = 487.75
$ws.Range("I135").Value = 496.26315
$ws.Range("J135").Value = 326
$ws.Range("K135").Value = 4466.36835
$ws.Range("L135").Value = 2934
$ws.Range("M135").Value = -1931.36835
$ws.Range("N135").Value = -8004
$ws.Range("H137").Value = 100002590
$ws.Range("I137").Value = 250001250
$ws.Range("K137").Value = 750003750
$ws.Range("M137").Value = -750001200
$ws.Range("H138").Value = 3463.7925
$ws.Range("I138").Value = 1484.2354
$ws.Range("J138").Value = 4398.5835
$ws.Range("K138").Value = 4452.706200000001
$ws.Range("L138").Value = 13195.7505
$ws.Range("M138").Value = 687.2937999999995
$ws.Range("N138").Value = -23475.7505

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2044.9722
$ws.Range("I2").Value = 2219.1482
$ws.Range("J2").Value = 1522.4445
$ws.Range("K2").Value = 2219.1482
$ws.Range("L2").Value = 1522.4445
$ws.Range("M2").Value = -2106.1482
$ws.Range("N2").Value = -1748.4445
$ws.Range("H25").Value = 2116.8
$ws.Range("I25").Value = 2194.6667
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 2194.6667
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = -1792.6667
$ws.Range("N25").Value = -2804
$ws.Range("H32").Value = 7577903.5
$ws.Range("I32").Value = 8929771
$ws.Range("K32").Value = 8929771
$ws.Range("M32").Value = -8929484
$ws.Range("H45").Value = 2594.7778
$ws.Range("I45").Value = 2800.9333
$ws.Range("J45").Value = 1564
$ws.Range("K45").Value = 2800.9333
$ws.Range("L45").Value = 1564
$ws.Range("M45").Value = -2423.9333
$ws.Range("N45").Value = -2318
$ws.Range("H61").Value = 2902.3809
$ws.Range("I61").Value = 2596.9678
$ws.Range("J61").Value = 3763.0908
$ws.Range("K61").Value = 2596.9678
$ws.Range("L61").Value = 3763.0908
$ws.Range("M61").Value = -2384.9678
$ws.Range("N61").Value = -4187.0908
$ws.Range("H74").Value = 2391.7585
$ws.Range("I74").Value = 1275.3889
$ws.Range("J74").Value = 4218.5454
$ws.Range("K74").Value = 1275.3889
$ws.Range("L74").Value = 4218.5454
$ws.Range("M74").Value = -401.3888999999999
$ws.Range("N74").Value = -5966.5454
$ws.Range("H77").Value = 2391.7585
$ws.Range("I77").Value = 1275.3889
$ws.Range("J77").Value = 4218.5454
$ws.Range("K77").Value = 6376.9445
$ws.Range("L77").Value = 21092.727
$ws.Range("M77").Value = -2008.9445
$ws.Range("N77").Value = -29828.727
$ws.Range("H97").Value = 1531.1482
$ws.Range("I97").Value = 1512.2727
$ws.Range("J97").Value = 1614.2
$ws.Range("K97").Value = 1512.2727
$ws.Range("L97").Value = 1614.2
$ws.Range("M97").Value = -1016.2727
$ws.Range("N97").Value = -2606.2
$ws.Range("H116").Value = 2044.9722
$ws.Range("I116").Value = 2219.1482
$ws.Range("J116").Value = 1522.4445
$ws.Range("K116").Value = 2219.1482
$ws.Range("L116").Value = 1522.4445
$ws.Range("M116").Value = 74.85179999999991
$ws.Range("N116").Value = -6110.4445
$ws.Range("H122").Value = 1084
$ws.Range("I122").Value = 811.3570999999999
$ws.Range("K122").Value = 2434.0713
$ws.Range("M122").Value = 15.92870000000039
$ws.Range("H132").Value = 1887.7797
$ws.Range("I132").Value = 1645.8148
$ws.Range("K132").Value = 4937.4444
$ws.Range("M132").Value = -2407.4444
$ws.Range("H136").Value = 2902.3809
$ws.Range("I136").Value = 2596.9678
$ws.Range("J136").Value = 3763.0908
$ws.Range("K136").Value = 7790.903399999999
$ws.Range("L136").Value = 11289.2724
$ws.Range("M136").Value = -5240.903399999999
$ws.Range("N136").Value = -16389.2724

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2044.9722
$ws.Range("I3").Value = 2219.1482
$ws.Range("J3").Value = 1522.4445
$ws.Range("K3").Value = 2219.1482
$ws.Range("L3").Value = 1522.4445
$ws.Range("M3").Value = -2105.1482
$ws.Range("N3").Value = -1750.4445
$ws.Range("H5").Value = 6407
$ws.Range("I5").Value = 442.75
$ws.Range("J5").Value = 12371.25
$ws.Range("K5").Value = 442.75
$ws.Range("L5").Value = 12371.25
$ws.Range("M5").Value = -329.75
$ws.Range("N5").Value = -12597.25
$ws.Range("H20").Value = 891.63635
$ws.Range("I20").Value = 840.9474
$ws.Range("J20").Value = 960.4286
$ws.Range("K20").Value = 840.9474
$ws.Range("L20").Value = 960.4286
$ws.Range("M20").Value = -593.9474
$ws.Range("N20").Value = -1454.4286
$ws.Range("H86").Value = 4071.1428
$ws.Range("I86").Value = 3866.6667
$ws.Range("J86").Value = 4224.5
$ws.Range("K86").Value = 3866.6667
$ws.Range("L86").Value = 4224.5
$ws.Range("M86").Value = -2743.6667
$ws.Range("N86").Value = -6470.5
$ws.Range("H89").Value = 4071.1428
$ws.Range("I89").Value = 3866.6667
$ws.Range("J89").Value = 4224.5
$ws.Range("K89").Value = 19333.3335
$ws.Range("L89").Value = 21122.5
$ws.Range("M89").Value = -13717.3335
$ws.Range("N89").Value = -32354.5
$ws.Range("H94").Value = 947.2222
$ws.Range("I94").Value = 967.5
$ws.Range("K94").Value = 967.5
$ws.Range("M94").Value = -516.5
$ws.Range("H134").Value = 23630.955
$ws.Range("I134").Value = 28448.777
$ws.Range("J134").Value = 4359.6665
$ws.Range("K134").Value = 85346.33099999999
$ws.Range("L134").Value = 13078.9995
$ws.Range("M134").Value = -82811.33099999999
$ws.Range("N134").Value = -18148.9995

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 101.695656
$ws.Range("I7").Value = 96.53333000000001
$ws.Range("J7").Value = 111.375
$ws.Range("K7").Value = 96.53333000000001
$ws.Range("L7").Value = 111.375
$ws.Range("M7").Value = 16.46666999999999
$ws.Range("N7").Value = -337.375
$ws.Range("H16").Value = 1647.2354
$ws.Range("I16").Value = 1519.1111
$ws.Range("J16").Value = 1791.375
$ws.Range("K16").Value = 1519.1111
$ws.Range("L16").Value = 1791.375
$ws.Range("M16").Value = -1232.1111
$ws.Range("N16").Value = -2365.375
$ws.Range("H31").Value = 88016.766
$ws.Range("I31").Value = 116461.43
$ws.Range("J31").Value = 8371.700000000001
$ws.Range("K31").Value = 116461.43
$ws.Range("L31").Value = 8371.700000000001
$ws.Range("M31").Value = -116166.43
$ws.Range("N31").Value = -8961.700000000001
$ws.Range("H34").Value = 88016.766
$ws.Range("I34").Value = 116461.43
$ws.Range("J34").Value = 8371.700000000001
$ws.Range("K34").Value = 116461.43
$ws.Range("L34").Value = 8371.700000000001
$ws.Range("M34").Value = -116259.43
$ws.Range("N34").Value = -8775.700000000001
$ws.Range("H58").Value = 2203.4722
$ws.Range("J58").Value = 3394.5334
$ws.Range("L58").Value = 3394.5334
$ws.Range("N58").Value = -3800.5334
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80630
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82184
$ws.Range("H80").Value = 50128
$ws.Range("J80").Value = 50128
$ws.Range("L80").Value = 50128
$ws.Range("N80").Value = -52374
$ws.Range("H83").Value = 50128
$ws.Range("J83").Value = 50128
$ws.Range("L83").Value = 150384
$ws.Range("N83").Value = -161616
$ws.Range("H92").Value = 79000.71000000001
$ws.Range("J92").Value = 80500.836
$ws.Range("L92").Value = 80500.836
$ws.Range("N92").Value = -85492.836
$ws.Range("H113").Value = 1647.2354
$ws.Range("I113").Value = 1519.1111
$ws.Range("J113").Value = 1791.375
$ws.Range("K113").Value = 1519.1111
$ws.Range("L113").Value = 1791.375
$ws.Range("M113").Value = 650.8888999999999
$ws.Range("N113").Value = -6131.375
$ws.Range("H134").Value = 23134.033
$ws.Range("I134").Value = 27504.4
$ws.Range("K134").Value = 82513.20000000001
$ws.Range("M134").Value = -79978.20000000001
$ws.Range("H136").Value = 2203.4722
$ws.Range("J136").Value = 3394.5334
$ws.Range("L136").Value = 10183.6002
$ws.Range("N136").Value = -15283.6002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 979
$ws.Range("I5").Value = 912
$ws.Range("J5").Value = 1113
$ws.Range("K5").Value = 2736
$ws.Range("L5").Value = 3339
$ws.Range("M5").Value = -2624
$ws.Range("N5").Value = -3563
$ws.Range("H11").Value = 2000439.8
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 2500499.8
$ws.Range("K11").Value = 600
$ws.Range("L11").Value = 7501499.399999999
$ws.Range("M11").Value = -460
$ws.Range("N11").Value = -7501779.399999999
$ws.Range("H25").Value = 1031.5555
$ws.Range("I25").Value = 930
$ws.Range("J25").Value = 1112.8
$ws.Range("K25").Value = 2790
$ws.Range("L25").Value = 3338.4
$ws.Range("M25").Value = -2621
$ws.Range("N25").Value = -3676.4
$ws.Range("H30").Value = 1031.5555
$ws.Range("I30").Value = 930
$ws.Range("J30").Value = 1112.8
$ws.Range("K30").Value = 2790
$ws.Range("L30").Value = 3338.4
$ws.Range("M30").Value = -2688
$ws.Range("N30").Value = -3542.4
$ws.Range("H31").Value = 675
$ws.Range("I31").Value = 685.7143
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 2057.1429
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = -1769.1429
$ws.Range("N31").Value = -2376
$ws.Range("H61").Value = 213.4
$ws.Range("I61").Value = 153
$ws.Range("J61").Value = 228.5
$ws.Range("K61").Value = 459
$ws.Range("L61").Value = 685.5
$ws.Range("M61").Value = -244
$ws.Range("N61").Value = -1115.5
$ws.Range("H92").Value = 677.55554
$ws.Range("I92").Value = 659.6
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 1978.8
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -730.8000000000002
$ws.Range("N92").Value = -4596
$ws.Range("H93").Value = 5295.2856
$ws.Range("J93").Value = 5548.846
$ws.Range("L93").Value = 16646.538
$ws.Range("N93").Value = -20390.538
$ws.Range("H107").Value = 346.45456
$ws.Range("I107").Value = 93
$ws.Range("J107").Value = 402.77777
$ws.Range("K107").Value = 279
$ws.Range("L107").Value = 1208.33331
$ws.Range("M107").Value = 1641
$ws.Range("N107").Value = -5048.33331
$ws.Range("H129").Value = 973.7143
$ws.Range("J129").Value = 1513.3334
$ws.Range("L129").Value = 4540.0002
$ws.Range("N129").Value = -14540.0002
$ws.Range("H131").Value = 1458.7715
$ws.Range("J131").Value = 1458.7715
$ws.Range("L131").Value = 4376.3145
$ws.Range("N131").Value = -14456.3145
$ws.Range("H135").Value = 979
$ws.Range("I135").Value = 912
$ws.Range("J135").Value = 1113
$ws.Range("K135").Value = 8208
$ws.Range("L135").Value = 10017
$ws.Range("M135").Value = -5673
$ws.Range("N135").Value = -15087
$ws.Range("H140").Value = 2692
$ws.Range("I140").Value = 1864.3043
$ws.Range("K140").Value = 5592.9129
$ws.Range("M140").Value = -412.9129000000003
$ws.Range("H141").Value = 2667
$ws.Range("I141").Value = 2602.923
$ws.Range("K141").Value = 7808.768999999999
$ws.Range("M141").Value = -2628.768999999999

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7114.923
$ws.Range("I70").Value = 7832.8887
$ws.Range("J70").Value = 5499.5
$ws.Range("K70").Value = 7832.8887
$ws.Range("L70").Value = 5499.5
$ws.Range("M70").Value = -7562.8887
$ws.Range("N70").Value = -6039.5
$ws.Range("H73").Value = 7114.923
$ws.Range("I73").Value = 7832.8887
$ws.Range("J73").Value = 5499.5
$ws.Range("K73").Value = 7832.8887
$ws.Range("L73").Value = 5499.5
$ws.Range("M73").Value = -6896.8887
$ws.Range("N73").Value = -7371.5
$ws.Range("H97").Value = 1368.1177
$ws.Range("I97").Value = 1824.7
$ws.Range("J97").Value = 715.8570999999999
$ws.Range("K97").Value = 1824.7
$ws.Range("L97").Value = 715.8570999999999
$ws.Range("M97").Value = -1328.7
$ws.Range("N97").Value = -1707.8571
$ws.Range("H113").Value = 2161.2
$ws.Range("I113").Value = 1753.5
$ws.Range("J113").Value = 2433
$ws.Range("K113").Value = 1753.5
$ws.Range("L113").Value = 2433
$ws.Range("M113").Value = 416.5
$ws.Range("N113").Value = -6773
$ws.Range("H118").Value = 49999
$ws.Range("J118").Value = 49999
$ws.Range("L118").Value = 49999
$ws.Range("N118").Value = -53313
$ws.Range("H119").Value = 80760.8
$ws.Range("J119").Value = 80760.8
$ws.Range("L119").Value = 80760.8
$ws.Range("N119").Value = -90436.8
$ws.Range("H122").Value = 67882.836
$ws.Range("I122").Value = 70459.60000000001
$ws.Range("K122").Value = 211378.8
$ws.Range("M122").Value = -208928.8
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -64900
$ws.Range("H132").Value = 23269304
$ws.Range("I132").Value = 32268966
$ws.Range("K132").Value = 96806898
$ws.Range("M132").Value = -96804368

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 901.7619
$ws.Range("I16").Value = 851.6429000000001
$ws.Range("J16").Value = 1002
$ws.Range("K16").Value = 851.6429000000001
$ws.Range("L16").Value = 1002
$ws.Range("M16").Value = -681.6429000000001
$ws.Range("N16").Value = -1342
$ws.Range("H22").Value = 1005.9
$ws.Range("I22").Value = 1038.625
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 1038.625
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -743.625
$ws.Range("N22").Value = -1465
$ws.Range("H27").Value = 1005.9
$ws.Range("I27").Value = 1038.625
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 1038.625
$ws.Range("L27").Value = 875
$ws.Range("M27").Value = -931.625
$ws.Range("N27").Value = -1089
$ws.Range("H40").Value = 3551.15
$ws.Range("I40").Value = 3240.0715
$ws.Range("K40").Value = 3240.0715
$ws.Range("M40").Value = -3104.0715
$ws.Range("H61").Value = 2470.9714
$ws.Range("J61").Value = 4899.3
$ws.Range("L61").Value = 4899.3
$ws.Range("N61").Value = -5303.3
$ws.Range("H68").Value = 4136.909
$ws.Range("I68").Value = 2380
$ws.Range("J68").Value = 5601
$ws.Range("K68").Value = 2380
$ws.Range("L68").Value = 5601
$ws.Range("M68").Value = -1631
$ws.Range("N68").Value = -7099
$ws.Range("H71").Value = 4136.909
$ws.Range("I71").Value = 2380
$ws.Range("J71").Value = 5601
$ws.Range("K71").Value = 11900
$ws.Range("L71").Value = 28005
$ws.Range("M71").Value = -8156
$ws.Range("N71").Value = -35493
$ws.Range("H74").Value = 59553.89
$ws.Range("I74").Value = 33300
$ws.Range("J74").Value = 72680.836
$ws.Range("K74").Value = 33300
$ws.Range("L74").Value = 72680.836
$ws.Range("M74").Value = -32302
$ws.Range("N74").Value = -74676.836
$ws.Range("H77").Value = 59553.89
$ws.Range("I77").Value = 33300
$ws.Range("J77").Value = 72680.836
$ws.Range("K77").Value = 99900
$ws.Range("L77").Value = 218042.508
$ws.Range("M77").Value = -94908
$ws.Range("N77").Value = -228026.508
$ws.Range("H93").Value = 2754.889
$ws.Range("I93").Value = 2326.8333
$ws.Range("J93").Value = 3611
$ws.Range("K93").Value = 2326.8333
$ws.Range("L93").Value = 3611
$ws.Range("M93").Value = -1078.8333
$ws.Range("N93").Value = -6107
$ws.Range("H108").Value = 96300.60000000001
$ws.Range("J108").Value = 96300.60000000001
$ws.Range("L108").Value = 96300.60000000001
$ws.Range("N108").Value = -103980.6
$ws.Range("H113").Value = 2470.9714
$ws.Range("J113").Value = 4899.3
$ws.Range("L113").Value = 4899.3
$ws.Range("N113").Value = -9239.299999999999
$ws.Range("H122").Value = 57871.668
$ws.Range("I122").Value = 1598.9231
$ws.Range("K122").Value = 4796.7693
$ws.Range("M122").Value = -2346.7693
$ws.Range("H132").Value = 2420.125
$ws.Range("I132").Value = 2322.3157
$ws.Range("J132").Value = 2791.8
$ws.Range("K132").Value = 6966.9471
$ws.Range("L132").Value = 8375.400000000001
$ws.Range("M132").Value = -4436.9471
$ws.Range("N132").Value = -13435.4
$ws.Range("H136").Value = 41539.688
$ws.Range("I136").Value = 1655.3846
$ws.Range("J136").Value = 214371.67
$ws.Range("K136").Value = 4966.1538
$ws.Range("L136").Value = 643115.01
$ws.Range("M136").Value = -2416.1538
$ws.Range("N136").Value = -648215.01

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 452399.75
$ws.Range("I62").Value = 899999.5
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 899999.5
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -899375.5
$ws.Range("N62").Value = -6048
$ws.Range("H64").Value = 57449.5
$ws.Range("J64").Value = 59899
$ws.Range("L64").Value = 59899
$ws.Range("N64").Value = -60395
$ws.Range("H65").Value = 452399.75
$ws.Range("I65").Value = 899999.5
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 4499997.5
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -4496877.5
$ws.Range("N65").Value = -30240
$ws.Range("H67").Value = 57449.5
$ws.Range("J67").Value = 59899
$ws.Range("L67").Value = 59899
$ws.Range("N67").Value = -61615
$ws.Range("H74").Value = 113134.37
$ws.Range("J74").Value = 120353.11
$ws.Range("L74").Value = 120353.11
$ws.Range("N74").Value = -122225.11
$ws.Range("H75").Value = 50020.75
$ws.Range("I75").Value = 53750
$ws.Range("J75").Value = 46291.5
$ws.Range("K75").Value = 53750
$ws.Range("L75").Value = 46291.5
$ws.Range("M75").Value = -52814
$ws.Range("N75").Value = -48163.5
$ws.Range("H77").Value = 113134.37
$ws.Range("J77").Value = 120353.11
$ws.Range("L77").Value = 361059.33
$ws.Range("N77").Value = -370419.33
$ws.Range("H78").Value = 50020.75
$ws.Range("I78").Value = 53750
$ws.Range("J78").Value = 46291.5
$ws.Range("K78").Value = 161250
$ws.Range("L78").Value = 138874.5
$ws.Range("M78").Value = -156570
$ws.Range("N78").Value = -148234.5
$ws.Range("H81").Value = 7094.8335
$ws.Range("I81").Value = 7094.8335
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 14189.667
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 7094.8335
$ws.Range("I84").Value = 7094.8335
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 70948.33499999999
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 2711.1538
$ws.Range("I107").Value = 1639.5625
$ws.Range("J107").Value = 4425.7
$ws.Range("K107").Value = 4918.6875
$ws.Range("L107").Value = 13277.1
$ws.Range("M107").Value = -2998.6875
$ws.Range("N107").Value = -17117.1
$ws.Range("H122").Value = 1307.4642
$ws.Range("I122").Value = 989.8421
$ws.Range("K122").Value = 2969.5263
$ws.Range("M122").Value = -519.5263
$ws.Range("H132").Value = 4547916.5
$ws.Range("I132").Value = 6251996
$ws.Range("K132").Value = 18755988
$ws.Range("M132").Value = -18753458
$ws.Range("H136").Value = 6438986
$ws.Range("I136").Value = 7778118
$ws.Range("K136").Value = 23334354
$ws.Range("M136").Value = -23331804
